$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "BKT0"
$ws.Range("B2").Value = "POPINDER SINGH"
$ws.Range("C2").Value = 817120112.4299997
$ws.Range("D2").Value = 228
$ws.Range("E2").Value = 25
$ws.Range("F2").Value = 202
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = 61842820.34
$ws.Range("N2").Value = 752412783.8699996
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = ""
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = 2864508.22
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = 7.57
$ws.Range("V2").Value = 92.08
$ws.Range("W2").Value = ""
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = 0.35
$ws.Range("Z2").Value = ""
$ws.Range("AA2").Value = ""
$ws.Range("AB2").Value = ""
$ws.Range("AC2").Value = 8915728
$ws.Range("AD2").Value = 92.42999999999999
$ws.Range("AE2").Value = 0.35

# Row 3
$ws.Range("A3").Value = "BKT1"
$ws.Range("B3").Value = "JOGINDER"
$ws.Range("C3").Value = 298414635.3600002
$ws.Range("D3").Value = 214
$ws.Range("E3").Value = 42
$ws.Range("F3").Value = 125
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = 44
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = 49824859.34
$ws.Range("N3").Value = 194592171.5200001
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = 50315521.08000001
$ws.Range("Q3").Value = ""
$ws.Range("R3").Value = ""
$ws.Range("S3").Value = 3682083.419999999
$ws.Range("T3").Value = ""
$ws.Range("U3").Value = 16.7
$ws.Range("V3").Value = 65.21
$ws.Range("W3").Value = ""
$ws.Range("X3").Value = ""
$ws.Range("Y3").Value = 1.23
$ws.Range("Z3").Value = 16.86
$ws.Range("AA3").Value = ""
$ws.Range("AB3").Value = ""
$ws.Range("AC3").Value = 5234507
$ws.Range("AD3").Value = 83.3
$ws.Range("AE3").Value = 18.09

# Row 4
$ws.Range("A4").Value = "BKT2"
$ws.Range("B4").Value = "JOGINDER"
$ws.Range("C4").Value = 1921240
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = 1921240
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = ""
$ws.Range("Q4").Value = ""
$ws.Range("R4").Value = ""
$ws.Range("S4").Value = ""
$ws.Range("T4").Value = ""
$ws.Range("U4").Value = 100
$ws.Range("V4").Value = ""
$ws.Range("W4").Value = ""
$ws.Range("X4").Value = ""
$ws.Range("Y4").Value = ""
$ws.Range("Z4").Value = ""
$ws.Range("AA4").Value = ""
$ws.Range("AB4").Value = ""
$ws.Range("AC4").Value = ""
$ws.Range("AD4").Value = ""
$ws.Range("AE4").Value = ""

# Row 5
$ws.Range("A5").Value = "BKT3"
$ws.Range("B5").Value = "JOGINDER"
$ws.Range("C5").Value = 1765526.64
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = 1765526.64
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = ""
$ws.Range("Q5").Value = ""
$ws.Range("R5").Value = ""
$ws.Range("S5").Value = ""
$ws.Range("T5").Value = ""
$ws.Range("U5").Value = 100
$ws.Range("V5").Value = ""
$ws.Range("W5").Value = ""
$ws.Range("X5").Value = ""
$ws.Range("Y5").Value = ""
$ws.Range("Z5").Value = ""
$ws.Range("AA5").Value = ""
$ws.Range("AB5").Value = ""
$ws.Range("AC5").Value = ""
$ws.Range("AD5").Value = ""
$ws.Range("AE5").Value = ""

# Row 6
$ws.Range("A6").Value = "BKT4"
$ws.Range("B6").Value = "JOGINDER"
$ws.Range("C6").Value = 489693
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = 489693
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = ""
$ws.Range("Q6").Value = ""
$ws.Range("R6").Value = ""
$ws.Range("S6").Value = ""
$ws.Range("T6").Value = ""
$ws.Range("U6").Value = 100
$ws.Range("V6").Value = ""
$ws.Range("W6").Value = ""
$ws.Range("X6").Value = ""
$ws.Range("Y6").Value = ""
$ws.Range("Z6").Value = ""
$ws.Range("AA6").Value = ""
$ws.Range("AB6").Value = ""
$ws.Range("AC6").Value = ""
$ws.Range("AD6").Value = ""
$ws.Range("AE6").Value = ""

# Row 7
$ws.Range("A7").Value = "BKT5"
$ws.Range("B7").Value = "JOGINDER"
$ws.Range("C7").Value = 2298733
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = 2298733
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = ""
$ws.Range("P7").Value = ""
$ws.Range("Q7").Value = ""
$ws.Range("R7").Value = ""
$ws.Range("S7").Value = ""
$ws.Range("T7").Value = ""
$ws.Range("U7").Value = 100
$ws.Range("V7").Value = ""
$ws.Range("W7").Value = ""
$ws.Range("X7").Value = ""
$ws.Range("Y7").Value = ""
$ws.Range("Z7").Value = ""
$ws.Range("AA7").Value = ""
$ws.Range("AB7").Value = ""
$ws.Range("AC7").Value = ""
$ws.Range("AD7").Value = ""
$ws.Range("AE7").Value = ""

# Row 8
$ws.Range("A8").Value = "BKT7"
$ws.Range("B8").Value = "JOGINDER"
$ws.Range("C8").Value = 1180938.14
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = 1180938.14
$ws.Range("N8").Value = ""
$ws.Range("O8").Value = ""
$ws.Range("P8").Value = ""
$ws.Range("Q8").Value = ""
$ws.Range("R8").Value = ""
$ws.Range("S8").Value = ""
$ws.Range("T8").Value = ""
$ws.Range("U8").Value = 100
$ws.Range("V8").Value = ""
$ws.Range("W8").Value = ""
$ws.Range("X8").Value = ""
$ws.Range("Y8").Value = ""
$ws.Range("Z8").Value = ""
$ws.Range("AA8").Value = ""
$ws.Range("AB8").Value = ""
$ws.Range("AC8").Value = ""
$ws.Range("AD8").Value = ""
$ws.Range("AE8").Value = ""

# Update used range / dimension implicitly handled by Excel; ensure rows below are set